# Update the dSF column (F) values as per the repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = 4
$ws.Range("F3").Value = 4
$ws.Range("F5").Value = 6
$ws.Range("F6").Value = -11
$ws.Range("F7").Value = -2
$ws.Range("F11").Value = 5
$ws.Range("F12").Value = 10
$ws.Range("F19").Value = 0
$ws.Range("F22").Value = 2
$ws.Range("F23").Value = -2

$wb.Save()
